# Weekly update: insert a new (more recent) weekly reading above the
# previous last row, shifting the old row 11 down to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11; existing row 11 (and below) shift down to row 12.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the latest weekly price data.
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44466
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101001
$ws.Range("J11").Value = "Arándano (blue)"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 160
$ws.Range("N11").Value = 13500
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 13750
$ws.Range("Q11").Value = "`$/bandeja 2 kilos"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 6875
$ws.Range("T11").Value = 2
